$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D3"  = "20.74"
    "D4"  = "6.219"
    "D6"  = "3.580"
    "D7"  = "6.572"
    "D8"  = "1.491"
    "D9"  = "0.8188"
    "D11" = "0.1630"
    "D12" = "0.08393"
    "D14" = "0.03210"
    "D15" = "0.09136"
    "D16" = "3.706"
    "D17" = "0.001639"
    "D18" = "0.04714"
    "D19" = "0.006510"
    "D20" = "0.006170"
    "D23" = "3.777"
    "D25" = "0.3355"
    "D26" = "0.1250"
    "D40" = "0.04696"
    "D41" = "0.007193"
    "D44" = "0.01110"
    "D45" = "0.00006602"
    "D46" = "0.00000000751"
    "D47" = "0.8106"
    "D48" = "0.002863"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.NumberFormat = "General"
}
